$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: refresh timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 07:22"

# --- Nueva Zelanda (row 71): "Casos criticos" (F) updated ---
$ws.Range("F71").Value = 3

# --- Guatemala / Sri Lanka swap places in the sorted list (rows 116-117) ---
# Row 116 becomes Sri Lanka with updated case counts
$ws.Range("A116").Value = "Sri Lanka"
$ws.Range("B116").Value = 295
$ws.Range("C116").Value = 24
$ws.Range("D116").Value = 96
$ws.Range("E116").Value = 192
$ws.Range("F116").Value = 1

# Row 117 becomes Guatemala with its updated case counts
$ws.Range("A117").Value = "Guatemala"
$ws.Range("B117").Value = 289
$ws.Range("C117").Value = 32
$ws.Range("D117").Value = 21
$ws.Range("E117").Value = 261
$ws.Range("F117").Value = 3

# --- Belice / Nueva Caledonia / Fiyi reorder (rows 182-184) ---
# Row 182 becomes Fiyi with updated case counts
$ws.Range("A182").Value = "Fiyi"
$ws.Range("C182").Value = 1
$ws.Range("D182").Value = 3
$ws.Range("E182").Value = 15
$ws.Range("F182").Value = 0
$ws.Range("H182").Value = 0

# Row 183 becomes Belice with its updated case counts
$ws.Range("A183").Value = "Belice"
$ws.Range("D183").Value = 2
$ws.Range("E183").Value = 14
$ws.Range("H183").Value = 2

# Row 184 becomes Nueva Caledonia with its updated case counts
$ws.Range("A184").Value = "Nueva Caledonia"
$ws.Range("B184").Value = 18
$ws.Range("D184").Value = 15
$ws.Range("E184").Value = 3
$ws.Range("F184").Value = 1
